$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 0.01779166666666667
$ws.Cells.Item(2, 8).Value = 0.053375
$ws.Cells.Item(2, 9).Value = 0.03510503888026929
$ws.Cells.Item(2, 10).Value = 0.03510503888026929
$ws.Cells.Item(2, 13).Value = 0.4047206666666667
$ws.Cells.Item(2, 14).Value = 1.214162
$ws.Cells.Item(2, 15).Value = 0.02578034976888792
$ws.Cells.Item(2, 16).Value = 0.02578034976888792
$ws.Cells.Item(2, 17).Value = 0.007200655194444445
$ws.Cells.Item(2, 18).Value = 0.06480589674999999
$ws.Cells.Item(2, 19).Value = 0.0009050201809837519
$ws.Cells.Item(2, 20).Value = 0.0009050201809837517

# Row 3
$ws.Cells.Item(3, 7).Value = 0.01779166666666667
$ws.Cells.Item(3, 8).Value = 0.053375
$ws.Cells.Item(3, 9).Value = 0.03510503888026929
$ws.Cells.Item(3, 10).Value = 0.03510503888026929
$ws.Cells.Item(3, 15).Value = 0.03222381288358415
$ws.Cells.Item(3, 16).Value = 0.03222381288358415
$ws.Cells.Item(3, 17).Value = 0.009000365305555556
$ws.Cells.Item(3, 18).Value = 0.08100328775
$ws.Cells.Item(3, 19).Value = 0.001131218204148744
$ws.Cells.Item(3, 20).Value = 0.001131218204148744

# Row 4
$ws.Cells.Item(4, 7).Value = 0.01779166666666667
$ws.Cells.Item(4, 8).Value = 0.053375
$ws.Cells.Item(4, 9).Value = 0.03510503888026929
$ws.Cells.Item(4, 10).Value = 0.03510503888026929
$ws.Cells.Item(4, 13).Value = 1.039987
$ws.Cells.Item(4, 14).Value = 3.119961
$ws.Cells.Item(4, 15).Value = 0.06624625531460326
$ws.Cells.Item(4, 16).Value = 0.06624625531460326
$ws.Cells.Item(4, 17).Value = 0.01850310204166667
$ws.Cells.Item(4, 18).Value = 0.166527918375
$ws.Cells.Item(4, 19).Value = 0.002325577368491394
$ws.Cells.Item(4, 20).Value = 0.002325577368491393

# Row 5
$ws.Cells.Item(5, 7).Value = 0.01779166666666667
$ws.Cells.Item(5, 8).Value = 0.053375
$ws.Cells.Item(5, 9).Value = 0.03510503888026929
$ws.Cells.Item(5, 10).Value = 0.03510503888026929
$ws.Cells.Item(5, 13).Value = 13.74822133333333
$ws.Cells.Item(5, 14).Value = 41.244664
$ws.Cells.Item(5, 15).Value = 0.8757495820329246
$ws.Cells.Item(5, 16).Value = 0.8757495820329247
$ws.Cells.Item(5, 17).Value = 0.2446037712222222
$ws.Cells.Item(5, 18).Value = 2.201433941
$ws.Cells.Item(5, 19).Value = 0.0307432231266454
$ws.Cells.Item(5, 20).Value = 0.0307432231266454

# Row 6
$ws.Cells.Item(6, 9).Value = 0.07168136529168917
$ws.Cells.Item(6, 10).Value = 0.07168136529168917
$ws.Cells.Item(6, 13).Value = 0.4047206666666667
$ws.Cells.Item(6, 14).Value = 1.214162
$ws.Cells.Item(6, 15).Value = 0.02578034976888792
$ws.Cells.Item(6, 16).Value = 0.02578034976888792
$ws.Cells.Item(6, 17).Value = 0.01470309709933333
$ws.Cells.Item(6, 18).Value = 0.132327873894
$ws.Cells.Item(6, 19).Value = 0.00184797066913117
$ws.Cells.Item(6, 20).Value = 0.00184797066913117

# Row 7
$ws.Cells.Item(7, 9).Value = 0.07168136529168917
$ws.Cells.Item(7, 10).Value = 0.07168136529168917
$ws.Cells.Item(7, 15).Value = 0.03222381288358415
$ws.Cells.Item(7, 16).Value = 0.03222381288358415
$ws.Cells.Item(7, 19).Value = 0.002309846902399235
$ws.Cells.Item(7, 20).Value = 0.002309846902399236

# Row 8
$ws.Cells.Item(8, 9).Value = 0.07168136529168917
$ws.Cells.Item(8, 10).Value = 0.07168136529168917
$ws.Cells.Item(8, 13).Value = 1.039987
$ws.Cells.Item(8, 14).Value = 3.119961
$ws.Cells.Item(8, 15).Value = 0.06624625531460326
$ws.Cells.Item(8, 16).Value = 0.06624625531460326
$ws.Cells.Item(8, 17).Value = 0.037781687723
$ws.Cells.Item(8, 18).Value = 0.340035189507
$ws.Cells.Item(8, 19).Value = 0.004748622026412581
$ws.Cells.Item(8, 20).Value = 0.004748622026412581

# Row 9
$ws.Cells.Item(9, 9).Value = 0.07168136529168917
$ws.Cells.Item(9, 10).Value = 0.07168136529168917
$ws.Cells.Item(9, 13).Value = 13.74822133333333
$ws.Cells.Item(9, 14).Value = 41.244664
$ws.Cells.Item(9, 15).Value = 0.8757495820329246
$ws.Cells.Item(9, 16).Value = 0.8757495820329247
$ws.Cells.Item(9, 17).Value = 0.4994591328186667
$ws.Cells.Item(9, 18).Value = 4.495132195368
$ws.Cells.Item(9, 19).Value = 0.06277492569374618
$ws.Cells.Item(9, 20).Value = 0.06277492569374619

# Row 10
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.09574300000000001
$ws.Cells.Item(10, 8).Value = 0.287229
$ws.Cells.Item(10, 9).Value = 0.1889121351295713
$ws.Cells.Item(10, 10).Value = 0.1889121351295713
$ws.Cells.Item(10, 13).Value = 0.4047206666666667
$ws.Cells.Item(10, 14).Value = 1.214162
$ws.Cells.Item(10, 15).Value = 0.02578034976888792
$ws.Cells.Item(10, 16).Value = 0.02578034976888792
$ws.Cells.Item(10, 17).Value = 0.03874917078866667
$ws.Cells.Item(10, 18).Value = 0.348742537098
$ws.Cells.Item(10, 19).Value = 0.004870220919227767
$ws.Cells.Item(10, 20).Value = 0.004870220919227767

# Row 11
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.09574300000000001
$ws.Cells.Item(11, 8).Value = 0.287229
$ws.Cells.Item(11, 9).Value = 0.1889121351295713
$ws.Cells.Item(11, 10).Value = 0.1889121351295713
$ws.Cells.Item(11, 15).Value = 0.03222381288358415
$ws.Cells.Item(11, 16).Value = 0.03222381288358415
$ws.Cells.Item(11, 17).Value = 0.04843402203933334
$ws.Cells.Item(11, 18).Value = 0.4359061983540001
$ws.Cells.Item(11, 19).Value = 0.00608746929385367
$ws.Cells.Item(11, 20).Value = 0.00608746929385367

# Row 12
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.09574300000000001
$ws.Cells.Item(12, 8).Value = 0.287229
$ws.Cells.Item(12, 9).Value = 0.1889121351295713
$ws.Cells.Item(12, 10).Value = 0.1889121351295713
$ws.Cells.Item(12, 13).Value = 1.039987
$ws.Cells.Item(12, 14).Value = 3.119961
$ws.Cells.Item(12, 15).Value = 0.06624625531460326
$ws.Cells.Item(12, 16).Value = 0.06624625531460326
$ws.Cells.Item(12, 17).Value = 0.09957147534100001
$ws.Cells.Item(12, 18).Value = 0.8961432780690001
$ws.Cells.Item(12, 19).Value = 0.01251472153582042
$ws.Cells.Item(12, 20).Value = 0.01251472153582041

# Row 13
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.09574300000000001
$ws.Cells.Item(13, 8).Value = 0.287229
$ws.Cells.Item(13, 9).Value = 0.1889121351295713
$ws.Cells.Item(13, 10).Value = 0.1889121351295713
$ws.Cells.Item(13, 13).Value = 13.74822133333333
$ws.Cells.Item(13, 14).Value = 41.244664
$ws.Cells.Item(13, 15).Value = 0.8757495820329246
$ws.Cells.Item(13, 16).Value = 0.8757495820329247
$ws.Cells.Item(13, 17).Value = 1.316295955117333
$ws.Cells.Item(13, 18).Value = 11.846663596056
$ws.Cells.Item(13, 19).Value = 0.1654397233806695
$ws.Cells.Item(13, 20).Value = 0.1654397233806695

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.3569486666666666
$ws.Cells.Item(14, 8).Value = 1.070846
$ws.Cells.Item(14, 9).Value = 0.7043014606984702
$ws.Cells.Item(14, 10).Value = 0.7043014606984702
$ws.Cells.Item(14, 13).Value = 0.4047206666666667
$ws.Cells.Item(14, 14).Value = 1.214162
$ws.Cells.Item(14, 15).Value = 0.02578034976888792
$ws.Cells.Item(14, 16).Value = 0.02578034976888792
$ws.Cells.Item(14, 17).Value = 0.1444645023391111
$ws.Cells.Item(14, 18).Value = 1.300180521052
$ws.Cells.Item(14, 19).Value = 0.01815713799954523
$ws.Cells.Item(14, 20).Value = 0.01815713799954523

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.3569486666666666
$ws.Cells.Item(15, 8).Value = 1.070846
$ws.Cells.Item(15, 9).Value = 0.7043014606984702
$ws.Cells.Item(15, 10).Value = 0.7043014606984702
$ws.Cells.Item(15, 15).Value = 0.03222381288358415
$ws.Cells.Item(15, 16).Value = 0.03222381288358415
$ws.Cells.Item(15, 17).Value = 0.1805715257328889
$ws.Cells.Item(15, 18).Value = 1.625143731596
$ws.Cells.Item(15, 19).Value = 0.0226952784831825
$ws.Cells.Item(15, 20).Value = 0.0226952784831825

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.3569486666666666
$ws.Cells.Item(16, 8).Value = 1.070846
$ws.Cells.Item(16, 9).Value = 0.7043014606984702
$ws.Cells.Item(16, 10).Value = 0.7043014606984702
$ws.Cells.Item(16, 13).Value = 1.039987
$ws.Cells.Item(16, 14).Value = 3.119961
$ws.Cells.Item(16, 15).Value = 0.06624625531460326
$ws.Cells.Item(16, 16).Value = 0.06624625531460326
$ws.Cells.Item(16, 17).Value = 0.3712219730006666
$ws.Cells.Item(16, 18).Value = 3.340997757006
$ws.Cells.Item(16, 19).Value = 0.04665733438387887
$ws.Cells.Item(16, 20).Value = 0.04665733438387887

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.3569486666666666
$ws.Cells.Item(17, 8).Value = 1.070846
$ws.Cells.Item(17, 9).Value = 0.7043014606984702
$ws.Cells.Item(17, 10).Value = 0.7043014606984702
$ws.Cells.Item(17, 13).Value = 13.74822133333333
$ws.Cells.Item(17, 14).Value = 41.244664
$ws.Cells.Item(17, 15).Value = 0.8757495820329246
$ws.Cells.Item(17, 16).Value = 0.8757495820329247
$ws.Cells.Item(17, 17).Value = 4.907409273971555
$ws.Cells.Item(17, 18).Value = 44.166683465744
$ws.Cells.Item(17, 19).Value = 0.6167917098318636
$ws.Cells.Item(17, 20).Value = 0.6167917098318637

Write-Output "Applied new TPM values"